$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "27.763.48"
$ws.Cells.Item(2, 5).Value = "  +3.11%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.868.29"
$ws.Cells.Item(3, 5).Value = "  +3.15%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.038"
$ws.Cells.Item(4, 5).Value = "  +2.99%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "324.44"
$ws.Cells.Item(5, 5).Value = "  +3.69%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "1.035"
$ws.Cells.Item(6, 5).Value = "  +2.82%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4420"
$ws.Cells.Item(7, 5).Value = "  +3.06%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3809"
$ws.Cells.Item(8, 5).Value = "  +3.04%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.07469"
$ws.Cells.Item(9, 5).Value = "  +3.03%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.8847"
$ws.Cells.Item(10, 5).Value = "  +2.52%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "21.80"
$ws.Cells.Item(11, 5).Value = "  +2.98%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "1.881.89"
$ws.Cells.Item(12, 5).Value = "  -8.15%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "5.567"
$ws.Cells.Item(13, 5).Value = "  +3.15%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "6.754"
$ws.Cells.Item(14, 5).Value = "  +1.71%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.07232"
$ws.Cells.Item(15, 5).Value = "  +4.46%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "83.80"
$ws.Cells.Item(16, 5).Value = "  +3.74%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "1.041"
$ws.Cells.Item(17, 5).Value = "  +3.34%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.000009115"
$ws.Cells.Item(18, 5).Value = "  +2.92%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "1.035"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "15.56"
$ws.Cells.Item(20, 5).Value = "  +2.39%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "27.797.55"
$ws.Cells.Item(21, 5).Value = "  +3.07%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "5.318"
$ws.Cells.Item(22, 5).Value = "  +2.24%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "11.46"
$ws.Cells.Item(23, 5).Value = "  +4.66%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "158.90"
$ws.Cells.Item(24, 5).Value = "  +3.09%  "
$ws.Cells.Item(25, 5).Value = "  +2.65%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "18.88"
$ws.Cells.Item(26, 5).Value = "  +2.82%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "1.994"
$ws.Cells.Item(27, 5).Value = "  +4.52%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "5.343"
$ws.Cells.Item(28, 5).Value = "  +2.05%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "117.78"
$ws.Cells.Item(29, 5).Value = "  +2.38%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.09108"
$ws.Cells.Item(30, 5).Value = "  +1.70%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.215"
$ws.Cells.Item(31, 5).Value = "  +5.09%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.7679"
$ws.Cells.Item(32, 5).Value = "  +3.37%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "4.580"
$ws.Cells.Item(33, 5).Value = "  +3.55%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "2.905"
$ws.Cells.Item(34, 5).Value = "  +3.53%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.036"
$ws.Cells.Item(35, 5).Value = "  +2.86%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "1.161"
$ws.Cells.Item(36, 5).Value = "  +3.35%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.01992"
$ws.Cells.Item(37, 5).Value = "  +3.22%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.05353"
$ws.Cells.Item(38, 5).Value = "  +2.35%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.5207"
$ws.Cells.Item(39, 5).Value = "  +2.35%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.845"
$ws.Cells.Item(40, 5).Value = "  +3.13%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.1696"
$ws.Cells.Item(41, 5).Value = "  +3.04%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "6.861"
$ws.Cells.Item(42, 5).Value = "  +6.34%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "8.712"
$ws.Cells.Item(43, 5).Value = "  +5.13%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "109.98"
$ws.Cells.Item(44, 5).Value = "  +2.63%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "10.67"
$ws.Cells.Item(45, 5).Value = "  +2.76%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "1.730"
$ws.Cells.Item(46, 5).Value = "  +4.81%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.4694"
$ws.Cells.Item(47, 5).Value = "  +3.00%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.06432"
$ws.Cells.Item(48, 5).Value = "  +2.40%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "1.857"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "39.87"
$ws.Cells.Item(50, 5).Value = "  +5.04%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.9369"
$ws.Cells.Item(51, 5).Value = "  +2.93%  "
